$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Replace the body text of the final (italic) paragraph while the
#    search string is still unique in the document (before we add a
#    second copy of it via the new "Meta description" paragraph).
# ------------------------------------------------------------------
$n = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($n)
$oldClosing = "Experience the underwater world of Atlantis. Play Aquaman slot online for free and enjoy unique mechanics, graphics, and four progressive jackpots."
$newClosing = "Create a feature image that perfectly captures the adventurous and fun spirit of the Aquaman slot game! The image should be in a cartoon style and prominently feature a happy Maya warrior with glasses. Make sure the Maya warrior is engaged in an exciting activity related to the game, such as spinning the reels with Aquaman or finding treasures in Atlantis. Use vibrant and bold colors to make the image standout, and include some of the game symbols like the trident, Mera, or the Aquaman symbol to tie everything together. Let the image showcase the thrill and excitement of this exciting game and make it the perfect visual representation of the Aquaman slot game."
$lastPara.Range.Find.Execute($oldClosing, $true, $false, $false, $false, $false, $true, 1, $false, $newClosing, 2) | Out-Null

# ------------------------------------------------------------------
# 2) Remove the trailing duplicate heading paragraph ("Play Aquaman
#    Slot for Free - Review & Demo 2021", bold) that now sits right
#    before the paragraph we just edited.
# ------------------------------------------------------------------
$n = $d.Paragraphs.Count
$dupHeading = $d.Paragraphs.Item($n - 1)
$dupHeading.Range.Delete() | Out-Null

# ------------------------------------------------------------------
# 3) Insert a new "Meta description" paragraph right after the title
#    (Heading1) paragraph at the top of the document.
# ------------------------------------------------------------------
$titlePara = $d.Paragraphs.Item(1)
$titlePara.Range.InsertParagraphAfter() | Out-Null

$metaPara = $d.Paragraphs.Item(2)
$metaRange = $metaPara.Range

$metaParagraphXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Meta description</w:t></w:r><w:r><w:t>: Experience the underwater world of Atlantis. Play Aquaman slot online for free and enjoy unique mechanics, graphics, and four progressive jackpots.</w:t></w:r></w:p>'
$metaRange.InsertXML($metaParagraphXml) | Out-Null

Write-Host "done"
